$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.342.66'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.550.06'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.08'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.480'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0889'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '1.772.12'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = '1.552.69'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('D14').Value = '28.322.97'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('E15').Value = '  -1.92%  '
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.76'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '227.49'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.91'
$ws.Range('E23').Value = '  -2.56%  '
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.11'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -1.96%  '
$ws.Range('E27').Value = '  -1.20%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('E29').Value = '  -3.10%  '
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('E31').Value = '  -4.74%  '
$ws.Range('D33').Value = '1.384.51'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.48'
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('E39').Value = '  -2.59%  '
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.510'
$ws.Range('E41').Value = '  -2.30%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.776'
$ws.Range('E43').Value = '  -2.19%  '
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.37'
$ws.Range('E45').Value = '  -1.56%  '
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('D47').Value = '1.685.15'
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.864'
$ws.Range('E48').Value = '  -10.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.42'
$ws.Range('E49').Value = '  -1.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.15'
$ws.Range('E50').Value = '  +5.11%  '
$ws.Range('E51').Value = '  -0.05%  '
